# Add a new option "Anamnese - Anotações Clínicas" to the sorted list of
# actions available for the professional, keeping the single column sorted
# alphabetically (as reflected by the worksheet's sortState) and updating
# the current selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Insert a new row right after the first one so every following item
# shifts down by one (old A2:A22 -> A3:A23) and place the new value there.
$ws.Rows("2:2").Insert()
$ws.Range("A2").Value = "Anamnese - Anotações Clínicas"

# Re-apply the alphabetical sort over A2:A23 (the first row is kept out of
# the sort, exactly like the original sortState/sortCondition on A1).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($ws.Range("A2:A23"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Make sure the final order exactly matches the expected alphabetical
# ordering of all the action labels.
$values = @(
  "ANTP . P",
  "Anamnese - Anotações Clínicas",
  "Atestados / Declarações",
  "Biomicroscopia",
  "C . R",
  "Concluir Atendimento",
  "Condutas e Justificativas",
  "Evolução Refrativa",
  "Fundo de Olho",
  "Fundo de Olho - Biomicroscopia",
  "Histórico de Solicitações",
  "Imagem(ns) de Exame(s)",
  "Laudar",
  "Mapeamento de Retina",
  "OCT - Tomografia de Coerência Óptica - Retina",
  "Portifólio",
  "Pressão Intra Ocular",
  "Receituário",
  "Retinografia",
  "Retorno",
  "Sair do Atendimento",
  "Solicitação de Procedimentos",
  "Solicitar Dilatação"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Update the active selection to the newly inserted cell.
[void]$ws.Range("A2").Select()
